# Locate the "Application Link" / video-URL content placeholder.
$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -like "*Application Link*") {
                $targetSlide = $sl
                $targetShape = $sh
            }
        }
    }
}

$tf = $targetShape.TextFrame
$tr = $tf.TextRange

# Find the paragraph that holds the "Application Link" text and the paragraph
# that holds the raw video URL (they are not adjacent - an empty paragraph
# sits between them). Split on the paragraph-mark carriage return to get an
# accurate paragraph count (TextRange.Paragraphs().Count can under-report
# trailing empty paragraphs).
$paraCount = ($tr.Text -split "`r").Count
$linkParaIdx = 0
$urlParaIdx = 0
for ($i = 1; $i -le $paraCount; $i++) {
    $txt = $tr.Paragraphs($i, 1).Text
    if ($txt -like "*Application Link*") {
        $linkParaIdx = $i
    } elseif ($txt -like "http*") {
        $urlParaIdx = $i
    }
}

$urlPara = $tr.Paragraphs($urlParaIdx, 1)
$videoUrl = $urlPara.Text

# Build the replacement text in-place on the URL paragraph, since its run
# formatting (no "dirty" flag) is the one we want to keep: prefix it with
# "Application Link " and then turn the URL portion into the word "Video".
$prefix = "Application Link "
$urlPara.Text = $prefix + $urlPara.Text

$urlPara = $tr.Paragraphs($urlParaIdx, 1)
$urlRange = $tr.Characters($urlPara.Start + $prefix.Length, $videoUrl.Length)
$urlRange.Text = "Video"

# Re-acquire the "Video" run and turn it into a hyperlink pointing at the
# original video URL.
$urlPara = $tr.Paragraphs($urlParaIdx, 1)
$videoRange = $tr.Characters($urlPara.Start + $prefix.Length, 5)
$actionSetting = $videoRange.ActionSettings(1)
$actionSetting.Hyperlink.Address = $videoUrl

# Remove the now-redundant leading paragraph(s): the original
# "Application Link" paragraph and the blank paragraph that separated it
# from the URL paragraph. Deleting the first paragraph repeatedly merges
# it away until only the (now combined) text paragraph and the original
# trailing blank paragraph remain.
for ($i = 1; $i -lt $urlParaIdx; $i++) {
    $tr.Paragraphs(1, 1).Delete()
}
